$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column J to column K for rows 3-6
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats

# Add year 2023 header in column K, row 3
$ws.Range("K3").Value = 2023

# Add new data values for 2023 in column K, rows 4-6
$ws.Range("K4").Value = 981.7
$ws.Range("K5").Value = 587.1
$ws.Range("K6").Value = 1324.2
